$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the tiny floating-point rounding differences in row 13
# (results of re-running the Gaussian Quadrature averaging)
$ws.Range("D13").Value = 0.9925452757136934
$ws.Range("J13").Value = 0.9925452757136934
$ws.Range("K13").Value = 0.9923159133273418
$ws.Range("L13").Value = 0.9918523608564576

# Add new row 16 for the "HexGrid-60degTilt5degRes" dataset (14th entry)
$ws.Range("A15").Copy($ws.Range("A16"))
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 1.365783696756323
$ws.Range("D16").Value = 0.9903387413958611
$ws.Range("E16").Value = 0.9629890908336594
$ws.Range("F16").Value = 1.365783696756323
$ws.Range("G16").Value = 0.9280892265155319
$ws.Range("H16").Value = 1.027240800069052
$ws.Range("I16").Value = 0.9858246838102366
$ws.Range("J16").Value = 0.9903387413958611
$ws.Range("K16").Value = 0.9766639161147602
$ws.Range("L16").Value = 1.171223806435542
$ws.Range("M16").Value = 1.043377706563444
